# Updates cryptos list values (prices and 1h volume change %) per source data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.143.72'
$ws.Range('E2').Value = '  -1.58%  '
$ws.Range('D3').Value = '2.270.23'
$ws.Range('E3').Value = '  -1.93%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '110.69'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.92%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '264.94'
$ws.Range('D6').Style = "Normal"
$ws.Range('E7').Value = '  -1.50%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -4.17%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '47.38'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.45%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0927'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.00%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '8.83'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.90%  '
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '15.40'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.91%  '
$ws.Range('D15').Value = '2.613.55'
$ws.Range('E15').Value = '  -1.77%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.851'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.84%  '
$ws.Range('D17').Value = '2.269.68'
$ws.Range('E17').Value = '  -1.40%  '
$ws.Range('D18').Value = '43.054.43'
$ws.Range('E18').Value = '  -1.79%  '
$ws.Range('E19').Value = '  -2.52%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.82'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.61%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '71.14'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.13%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.47'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.88%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '231.25'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.55%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.64'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.34%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.86'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.45%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.29'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.54%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '3.91'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.32%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '40.26'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -7.49%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.25'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.94%  '
$ws.Range('B31').Value = 'WEMIXToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.30'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -4.07%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '171.87'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -3.43%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '21.25'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -2.93%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0905'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -3.58%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.78'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +2.27%  '
$ws.Range('E36').Value = '  -0.61%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.69'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.65%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0353'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.83%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.83'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -3.40%  '
$ws.Range('E40').Value = '  -6.75%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.64'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +9.51%  '
$ws.Range('B42').Value = 'MultiversX'
$ws.Range('C42').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '74.87'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +7.89%  '
$ws.Range('B43').Value = 'Celestia'
$ws.Range('C43').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '13.93'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +9.89%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.236'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -3.99%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '6.11'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +6.77%  '
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('E47').Value = '  -0.19%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '8.65'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.35%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0992'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.09%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.25'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.52%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '100.78'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +0.39%  '
